# Add Week 05 & Week 06 deliverables:
# a new "Estimated Cost (RM)" column (H) with per-task cost figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 --------------------------------------------------
$ws.Range("H1").Value = "Estimated Cost (RM)"
# Match the look of the existing header row (bold, centered, wrapped).
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H1").VerticalAlignment = -4108     # xlCenter
$ws.Range("H1").WrapText = $true

# Header row needs to grow to fit the new, taller wrapped header text.
$ws.Rows.Item(1).RowHeight = 43.2

# --- Estimated cost values for each task (rows 2-9) -------------------
$costs = @{
    2 = 100
    3 = 200
    4 = 200
    5 = 400
    6 = 100
    7 = 100
    8 = 200
    9 = 300
}

# Match the look of the existing body rows (vertical-centered, wrapped) --
# row 2 is left at the default style, matching the rest of the column.
$ws.Range("H3:H9").VerticalAlignment = -4108  # xlCenter
$ws.Range("H3:H9").WrapText = $true

foreach ($row in $costs.Keys) {
    $ws.Cells.Item($row, 8).Value = $costs[$row]
}

# Leave the selection where it naturally lands after filling the column.
$ws.Range("I9").Select() | Out-Null
